$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Header text updates (Volume/Number and the reporting week date range).
# We edit via Characters(start, length).Text so only the digits/dates are
# touched and the surrounding text ("Volume ", "   Number  ", etc.) as
# well as the run formatting is left intact.
# -----------------------------------------------------------------------

# A8 merged cell (A8:B8): "Volume 32   Number  50" -> "...51"
$a8 = $ws.Cells.Item(8, 1)
$a8full = $a8.Value2
$numStart = $a8full.Length - 1   # "50" is the last two characters (1-based start)
$a8.Characters($numStart, 2).Text = "51"

# C9 merged cell (C9:L9): week range 12/8/2025 -> 12/14/2025 becomes 12/15/2025 -> 12/21/2025
# Replace the later date first so the earlier date's character offset is unaffected
# by the length change of the first replacement (9 chars -> 10 chars).
$c9 = $ws.Cells.Item(9, 3)
$c9full = $c9.Value2
$endIdx = $c9full.IndexOf("12/14/2025") + 1   # 1-based start position
$c9.Characters($endIdx, 10).Text = "12/21/2025"

$c9full2 = $c9.Value2
$startIdx = $c9full2.IndexOf("12/8/2025") + 1   # 1-based start position
$c9.Characters($startIdx, 9).Text = "12/15/2025"


# --- Reference cells used to copy cell formats (number format / style) when a cell
# changes between "text placeholder" (e.g. "0" / "***.*") and numeric representation. ---
$srcText = $ws.Cells.Item(23, 3)   # C23 -> text style (General / numFmt 0)
$srcInt  = $ws.Cells.Item(16, 3)   # C16 -> integer style (#,##0 like)
$srcPct  = $ws.Cells.Item(16, 8)   # H16 -> percent-change style (#,##0.0 like)

# ===== Row 14 =====
# F14: convert text-placeholder cell to integer number 1
$cell = $ws.Cells.Item(14, 6)
$srcInt.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 1

# H14: -> 0
$cell = $ws.Cells.Item(14, 8)
$cell.Value = 0

# I14: convert text-placeholder cell to integer number 1
$cell = $ws.Cells.Item(14, 9)
$srcInt.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 1

# K14: -> -75
$cell = $ws.Cells.Item(14, 11)
$cell.Value = -75

# M14: -> -50
$cell = $ws.Cells.Item(14, 13)
$cell.Value = -50

# N14: -> -90.909090909090
$cell = $ws.Cells.Item(14, 14)
$cell.Value = -90.909090909090

# ===== Row 15 =====
# C15: -> 2
$cell = $ws.Cells.Item(15, 3)
$cell.Value = 2

# F15: -> 3
$cell = $ws.Cells.Item(15, 6)
$cell.Value = 3

# G15: convert numeric-style cell to text placeholder "0"
$cell = $ws.Cells.Item(15, 7)
$srcText.Copy()
$cell.PasteSpecial(-4122)
$cell.NumberFormat = "@"
$cell.Value = "0"
$cell.NumberFormat = "General"

# H15: convert numeric-style cell to text placeholder "***.*"
$cell = $ws.Cells.Item(15, 8)
$srcText.Copy()
$cell.PasteSpecial(-4122)
$cell.NumberFormat = "@"
$cell.Value = "***.*"
$cell.NumberFormat = "General"

# I15: -> 18
$cell = $ws.Cells.Item(15, 9)
$cell.Value = 18

# K15: -> 20
$cell = $ws.Cells.Item(15, 11)
$cell.Value = 20

# L15: -> 5.882352941176
$cell = $ws.Cells.Item(15, 12)
$cell.Value = 5.882352941176

# M15: -> 350
$cell = $ws.Cells.Item(15, 13)
$cell.Value = 350

# N15: -> 63.636363636363
$cell = $ws.Cells.Item(15, 14)
$cell.Value = 63.636363636363

# ===== Row 16 =====
# D16: -> 7
$cell = $ws.Cells.Item(16, 4)
$cell.Value = 7

# E16: -> -71.428571428571
$cell = $ws.Cells.Item(16, 5)
$cell.Value = -71.428571428571

# G16: -> 13
$cell = $ws.Cells.Item(16, 7)
$cell.Value = 13

# H16: -> -38.461538461538
$cell = $ws.Cells.Item(16, 8)
$cell.Value = -38.461538461538

# I16: -> 150
$cell = $ws.Cells.Item(16, 9)
$cell.Value = 150

# J16: -> 147
$cell = $ws.Cells.Item(16, 10)
$cell.Value = 147

# K16: -> 2.040816326530
$cell = $ws.Cells.Item(16, 11)
$cell.Value = 2.040816326530

# L16: -> 16.279069767441
$cell = $ws.Cells.Item(16, 12)
$cell.Value = 16.279069767441

# M16: -> 59.574468085106
$cell = $ws.Cells.Item(16, 13)
$cell.Value = 59.574468085106

# N16: -> -82.248520710059
$cell = $ws.Cells.Item(16, 14)
$cell.Value = -82.248520710059

# ===== Row 17 =====
# C17: -> 2
$cell = $ws.Cells.Item(17, 3)
$cell.Value = 2

# D17: -> 7
$cell = $ws.Cells.Item(17, 4)
$cell.Value = 7

# E17: -> -71.428571428571
$cell = $ws.Cells.Item(17, 5)
$cell.Value = -71.428571428571

# G17: -> 14
$cell = $ws.Cells.Item(17, 7)
$cell.Value = 14

# H17: -> -7.142857142857
$cell = $ws.Cells.Item(17, 8)
$cell.Value = -7.142857142857

# I17: -> 155
$cell = $ws.Cells.Item(17, 9)
$cell.Value = 155

# J17: -> 193
$cell = $ws.Cells.Item(17, 10)
$cell.Value = 193

# K17: -> -19.689119170984
$cell = $ws.Cells.Item(17, 11)
$cell.Value = -19.689119170984

# L17: -> 8.391608391608
$cell = $ws.Cells.Item(17, 12)
$cell.Value = 8.391608391608

# M17: -> 127.941176470588
$cell = $ws.Cells.Item(17, 13)
$cell.Value = 127.941176470588

# N17: -> -12.921348314606
$cell = $ws.Cells.Item(17, 14)
$cell.Value = -12.921348314606

# ===== Row 18 =====
# C18: -> 5
$cell = $ws.Cells.Item(18, 3)
$cell.Value = 5

# D18: -> 2
$cell = $ws.Cells.Item(18, 4)
$cell.Value = 2

# E18: -> 150
$cell = $ws.Cells.Item(18, 5)
$cell.Value = 150

# F18: -> 17
$cell = $ws.Cells.Item(18, 6)
$cell.Value = 17

# G18: -> 12
$cell = $ws.Cells.Item(18, 7)
$cell.Value = 12

# H18: -> 41.666666666666
$cell = $ws.Cells.Item(18, 8)
$cell.Value = 41.666666666666

# I18: -> 196
$cell = $ws.Cells.Item(18, 9)
$cell.Value = 196

# J18: -> 194
$cell = $ws.Cells.Item(18, 10)
$cell.Value = 194

# K18: -> 1.030927835051
$cell = $ws.Cells.Item(18, 11)
$cell.Value = 1.030927835051

# L18: -> 0.512820512820
$cell = $ws.Cells.Item(18, 12)
$cell.Value = 0.512820512820

# M18: -> 16.666666666666
$cell = $ws.Cells.Item(18, 13)
$cell.Value = 16.666666666666

# N18: -> -77.777777777777
$cell = $ws.Cells.Item(18, 14)
$cell.Value = -77.777777777777

# ===== Row 19 =====
# C19: -> 21
$cell = $ws.Cells.Item(19, 3)
$cell.Value = 21

# D19: -> 35
$cell = $ws.Cells.Item(19, 4)
$cell.Value = 35

# E19: -> -40
$cell = $ws.Cells.Item(19, 5)
$cell.Value = -40

# F19: -> 106
$cell = $ws.Cells.Item(19, 6)
$cell.Value = 106

# G19: -> 104
$cell = $ws.Cells.Item(19, 7)
$cell.Value = 104

# H19: -> 1.923076923076
$cell = $ws.Cells.Item(19, 8)
$cell.Value = 1.923076923076

# I19: -> 1203
$cell = $ws.Cells.Item(19, 9)
$cell.Value = 1203

# J19: -> 1173
$cell = $ws.Cells.Item(19, 10)
$cell.Value = 1173

# K19: -> 2.557544757033
$cell = $ws.Cells.Item(19, 11)
$cell.Value = 2.557544757033

# L19: -> -6.381322957198
$cell = $ws.Cells.Item(19, 12)
$cell.Value = -6.381322957198

# M19: -> 15.229885057471
$cell = $ws.Cells.Item(19, 13)
$cell.Value = 15.229885057471

# N19: -> -68.275316455696
$cell = $ws.Cells.Item(19, 14)
$cell.Value = -68.275316455696

# ===== Row 20 =====
# D20: convert text-placeholder cell to integer number 1
$cell = $ws.Cells.Item(20, 4)
$srcInt.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 1

# E20: convert text-placeholder cell to percent-change number 0
$cell = $ws.Cells.Item(20, 5)
$srcPct.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 0

# F20: -> 3
$cell = $ws.Cells.Item(20, 6)
$cell.Value = 3

# G20: -> 3
$cell = $ws.Cells.Item(20, 7)
$cell.Value = 3

# I20: -> 35
$cell = $ws.Cells.Item(20, 9)
$cell.Value = 35

# J20: -> 44
$cell = $ws.Cells.Item(20, 10)
$cell.Value = 44

# K20: -> -20.454545454545
$cell = $ws.Cells.Item(20, 11)
$cell.Value = -20.454545454545

# L20: -> -54.545454545454
$cell = $ws.Cells.Item(20, 12)
$cell.Value = -54.545454545454

# M20: -> -7.894736842105
$cell = $ws.Cells.Item(20, 13)
$cell.Value = -7.894736842105

# N20: -> -96.013667425968
$cell = $ws.Cells.Item(20, 14)
$cell.Value = -96.013667425968

# ===== Row 21 =====
# C21: -> 33
$cell = $ws.Cells.Item(21, 3)
$cell.Value = 33

# D21: -> 52
$cell = $ws.Cells.Item(21, 4)
$cell.Value = 52

# E21: -> -36.538461538461
$cell = $ws.Cells.Item(21, 5)
$cell.Value = -36.538461538461

# F21: -> 151
$cell = $ws.Cells.Item(21, 6)
$cell.Value = 151

# G21: -> 147
$cell = $ws.Cells.Item(21, 7)
$cell.Value = 147

# H21: -> 2.721088435374
$cell = $ws.Cells.Item(21, 8)
$cell.Value = 2.721088435374

# I21: -> 1758
$cell = $ws.Cells.Item(21, 9)
$cell.Value = 1758

# J21: -> 1770
$cell = $ws.Cells.Item(21, 10)
$cell.Value = 1770

# K21: -> -0.677966101694
$cell = $ws.Cells.Item(21, 11)
$cell.Value = -0.677966101694

# L21: -> -4.767063921993
$cell = $ws.Cells.Item(21, 12)
$cell.Value = -4.767063921993

# M21: -> 23.977433004231
$cell = $ws.Cells.Item(21, 13)
$cell.Value = 23.977433004231

# N21: -> -73.351523419736
$cell = $ws.Cells.Item(21, 14)
$cell.Value = -73.351523419736

# ===== Row 22 =====
# D22: convert text-placeholder cell to integer number 3
$cell = $ws.Cells.Item(22, 4)
$srcInt.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 3

# E22: convert text-placeholder cell to percent-change number 0
$cell = $ws.Cells.Item(22, 5)
$srcPct.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 0

# F22: -> 11
$cell = $ws.Cells.Item(22, 6)
$cell.Value = 11

# G22: -> 3
$cell = $ws.Cells.Item(22, 7)
$cell.Value = 3

# H22: -> 266.666666666667
$cell = $ws.Cells.Item(22, 8)
$cell.Value = 266.666666666667

# I22: -> 100
$cell = $ws.Cells.Item(22, 9)
$cell.Value = 100

# J22: -> 91
$cell = $ws.Cells.Item(22, 10)
$cell.Value = 91

# K22: -> 9.890109890109
$cell = $ws.Cells.Item(22, 11)
$cell.Value = 9.890109890109

# L22: -> -9.909909909909
$cell = $ws.Cells.Item(22, 12)
$cell.Value = -9.909909909909

# M22: -> 44.927536231884
$cell = $ws.Cells.Item(22, 13)
$cell.Value = 44.927536231884

# ===== Row 24 =====
# D24: -> 82
$cell = $ws.Cells.Item(24, 4)
$cell.Value = 82

# E24: -> -6.097560975609
$cell = $ws.Cells.Item(24, 5)
$cell.Value = -6.097560975609

# F24: -> 281
$cell = $ws.Cells.Item(24, 6)
$cell.Value = 281

# G24: -> 295
$cell = $ws.Cells.Item(24, 7)
$cell.Value = 295

# H24: -> -4.745762711864
$cell = $ws.Cells.Item(24, 8)
$cell.Value = -4.745762711864

# I24: -> 3419
$cell = $ws.Cells.Item(24, 9)
$cell.Value = 3419

# J24: -> 3793
$cell = $ws.Cells.Item(24, 10)
$cell.Value = 3793

# K24: -> -9.860268916424
$cell = $ws.Cells.Item(24, 11)
$cell.Value = -9.860268916424

# L24: -> -16.28305582762
$cell = $ws.Cells.Item(24, 12)
$cell.Value = -16.28305582762

# M24: -> 106.336753168377
$cell = $ws.Cells.Item(24, 13)
$cell.Value = 106.336753168377

# ===== Row 25 =====
# C25: -> 69
$cell = $ws.Cells.Item(25, 3)
$cell.Value = 69

# D25: -> 81
$cell = $ws.Cells.Item(25, 4)
$cell.Value = 81

# E25: -> -14.814814814814
$cell = $ws.Cells.Item(25, 5)
$cell.Value = -14.814814814814

# F25: -> 267
$cell = $ws.Cells.Item(25, 6)
$cell.Value = 267

# G25: -> 296
$cell = $ws.Cells.Item(25, 7)
$cell.Value = 296

# H25: -> -9.797297297297
$cell = $ws.Cells.Item(25, 8)
$cell.Value = -9.797297297297

# I25: -> 3317
$cell = $ws.Cells.Item(25, 9)
$cell.Value = 3317

# J25: -> 3676
$cell = $ws.Cells.Item(25, 10)
$cell.Value = 3676

# K25: -> -9.766050054406
$cell = $ws.Cells.Item(25, 11)
$cell.Value = -9.766050054406

# L25: -> -18.800489596083
$cell = $ws.Cells.Item(25, 12)
$cell.Value = -18.800489596083

# ===== Row 26 =====
# C26: -> 4
$cell = $ws.Cells.Item(26, 3)
$cell.Value = 4

# D26: -> 12
$cell = $ws.Cells.Item(26, 4)
$cell.Value = 12

# E26: -> -66.666666666666
$cell = $ws.Cells.Item(26, 5)
$cell.Value = -66.666666666666

# G26: -> 32
$cell = $ws.Cells.Item(26, 7)
$cell.Value = 32

# H26: -> -9.375
$cell = $ws.Cells.Item(26, 8)
$cell.Value = -9.375

# I26: -> 434
$cell = $ws.Cells.Item(26, 9)
$cell.Value = 434

# J26: -> 442
$cell = $ws.Cells.Item(26, 10)
$cell.Value = 442

# K26: -> -1.809954751131
$cell = $ws.Cells.Item(26, 11)
$cell.Value = -1.809954751131

# L26: -> 7.425742574257
$cell = $ws.Cells.Item(26, 12)
$cell.Value = 7.425742574257

# M26: -> 64.393939393939
$cell = $ws.Cells.Item(26, 13)
$cell.Value = 64.393939393939

# ===== Row 27 =====
# C27: -> 2
$cell = $ws.Cells.Item(27, 3)
$cell.Value = 2

# F27: -> 3
$cell = $ws.Cells.Item(27, 6)
$cell.Value = 3

# G27: convert numeric-style cell to text placeholder "0"
$cell = $ws.Cells.Item(27, 7)
$srcText.Copy()
$cell.PasteSpecial(-4122)
$cell.NumberFormat = "@"
$cell.Value = "0"
$cell.NumberFormat = "General"

# H27: convert numeric-style cell to text placeholder "***.*"
$cell = $ws.Cells.Item(27, 8)
$srcText.Copy()
$cell.PasteSpecial(-4122)
$cell.NumberFormat = "@"
$cell.Value = "***.*"
$cell.NumberFormat = "General"

# I27: -> 21
$cell = $ws.Cells.Item(27, 9)
$cell.Value = 21

# K27: -> 10.526315789473
$cell = $ws.Cells.Item(27, 11)
$cell.Value = 10.526315789473

# L27: -> 5
$cell = $ws.Cells.Item(27, 12)
$cell.Value = 5

# ===== Row 28 =====
# C28: -> 1
$cell = $ws.Cells.Item(28, 3)
$cell.Value = 1

# D28: convert numeric-style cell to text placeholder "0"
$cell = $ws.Cells.Item(28, 4)
$srcText.Copy()
$cell.PasteSpecial(-4122)
$cell.NumberFormat = "@"
$cell.Value = "0"
$cell.NumberFormat = "General"

# E28: convert numeric-style cell to text placeholder "***.*"
$cell = $ws.Cells.Item(28, 5)
$srcText.Copy()
$cell.PasteSpecial(-4122)
$cell.NumberFormat = "@"
$cell.Value = "***.*"
$cell.NumberFormat = "General"

# F28: -> 5
$cell = $ws.Cells.Item(28, 6)
$cell.Value = 5

# G28: -> 3
$cell = $ws.Cells.Item(28, 7)
$cell.Value = 3

# H28: -> 66.666666666666
$cell = $ws.Cells.Item(28, 8)
$cell.Value = 66.666666666666

# I28: -> 98
$cell = $ws.Cells.Item(28, 9)
$cell.Value = 98

# K28: -> 3.157894736842
$cell = $ws.Cells.Item(28, 11)
$cell.Value = 3.157894736842

# L28: -> 4.255319148936
$cell = $ws.Cells.Item(28, 12)
$cell.Value = 4.255319148936

# ===== Row 31 =====
# D31: convert text-placeholder cell to integer number 1
$cell = $ws.Cells.Item(31, 4)
$srcInt.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 1

# E31: convert text-placeholder cell to percent-change number -100
$cell = $ws.Cells.Item(31, 5)
$srcPct.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = -100

# G31: convert text-placeholder cell to integer number 1
$cell = $ws.Cells.Item(31, 7)
$srcInt.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = 1

# H31: convert text-placeholder cell to percent-change number -100
$cell = $ws.Cells.Item(31, 8)
$srcPct.Copy()
$cell.PasteSpecial(-4122)
$cell.Value = -100

# J31: -> 14
$cell = $ws.Cells.Item(31, 10)
$cell.Value = 14

# K31: -> -14.285714285714
$cell = $ws.Cells.Item(31, 11)
$cell.Value = -14.285714285714
